# Apply crypto price/volume updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.539.87"
$ws.Range("E2").Value = "  -2.46%  "
$ws.Range("D3").Value = "1.978.94"
$ws.Range("E3").Value = "  -3.32%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "'244.44"
$ws.Range("E5").Value = "  +1.16%  "
$ws.Range("D6").Value = "'0.634"
$ws.Range("E6").Value = "  -4.82%  "
$ws.Range("D7").Value = "'57.44"
$ws.Range("E7").Value = "  +6.17%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'59.22"
$ws.Range("E9").Value = "  +1.76%  "
$ws.Range("D10").Value = "'0.359"
$ws.Range("E10").Value = "  +0.82%  "
$ws.Range("D11").Value = "'0.0730"
$ws.Range("E11").Value = "  -1.97%  "
$ws.Range("E12").Value = "  -2.91%  "
$ws.Range("D13").Value = "'0.931"
$ws.Range("E13").Value = "  +4.84%  "
$ws.Range("D14").Value = "'14.20"
$ws.Range("E14").Value = "  -3.31%  "
$ws.Range("D15").Value = "2.271.85"
$ws.Range("E15").Value = "  -3.22%  "
$ws.Range("D16").Value = "'5.23"
$ws.Range("E16").Value = "  -2.07%  "
$ws.Range("D17").Value = "1.977.93"
$ws.Range("E17").Value = "  -3.55%  "
$ws.Range("D18").Value = "'17.22"
$ws.Range("E18").Value = "  +4.67%  "
$ws.Range("D19").Value = "35.471.62"
$ws.Range("E19").Value = "  -2.58%  "
$ws.Range("D20").Value = "'70.86"
$ws.Range("E20").Value = "  -1.11%  "
$ws.Range("D21").Value = "0.0₃0844"
$ws.Range("E21").Value = "  -0.43%  "
$ws.Range("D22").Value = "'232.43"
$ws.Range("E22").Value = "  -1.68%  "
$ws.Range("D23").Value = "'5.14"
$ws.Range("E23").Value = "  -1.79%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "'2.49"
$ws.Range("E25").Value = "  +18.22%  "
$ws.Range("D26").Value = "'2.30"
$ws.Range("E26").Value = "  -1.70%  "
$ws.Range("D27").Value = "'9.13"
$ws.Range("E27").Value = "  -1.79%  "
$ws.Range("D28").Value = "'163.21"
$ws.Range("E28").Value = "  +0.59%  "
$ws.Range("D29").Value = "'19.24"
$ws.Range("E29").Value = "  -4.01%  "
$ws.Range("D30").Value = "'0.118"
$ws.Range("E30").Value = "  -2.23%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'4.83"
$ws.Range("E31").Value = "  -4.74%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'1.13"
$ws.Range("E32").Value = "  -0.90%  "
$ws.Range("D33").Value = "'0.0596"
$ws.Range("E33").Value = "  +1.58%  "
$ws.Range("D34").Value = "'0.0914"
$ws.Range("E34").Value = "  +12.17%  "
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").Value = "'4.25"
$ws.Range("E35").Value = "  -5.69%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "'2.34"
$ws.Range("E36").Value = "  +6.97%  "
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("D38").Value = "'1.76"
$ws.Range("E38").Value = "  -5.38%  "
$ws.Range("D39").Value = "'5.05"
$ws.Range("E39").Value = "  +5.95%  "
$ws.Range("E40").Value = "  -3.53%  "
$ws.Range("D41").Value = "'2.85"
$ws.Range("E41").Value = "  +2.38%  "
$ws.Range("D42").Value = "'0.0210"
$ws.Range("E42").Value = "  -1.26%  "
$ws.Range("E43").Value = "  -1.39%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "'16.04"
$ws.Range("E44").Value = "  +3.48%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.386.16"
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'91.11"
$ws.Range("E46").Value = "  -2.03%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.0881"
$ws.Range("E47").Value = "  -1.14%  "
$ws.Range("E48").Value = "  +2.59%  "
$ws.Range("E49").Value = "  +2.31%  "
$ws.Range("B50").Value = "FTXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D50").Value = "'3.76"
$ws.Range("E50").Value = "  +18.62%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "'2.27"
$ws.Range("E51").Value = "  -0.10%  "
